# Fruta / hortaliza, semanal
#
# The weekly price-report row that used to sit at row 141 (date 2021-02-11 /
# serial 44238) is being replaced by an updated weekly reading (date
# 2021-09-09 / serial 44448, with refreshed Volumen/Precio/S values), and the
# original 2021-02-11 row is kept as its own (new) row, pushing the old
# row 142 down to row 143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at 142 - this shifts the former row 142 (and
# anything below it) down to row 143, and extends the sheet dimension to
# A1:T143.
$ws.Rows.Item(142).Insert()

# Row 141 keeps its other fields, but gets a new date and updated
# Volumen / Precio minimo / Precio maximo / Precio promedio ponderado /
# Precio $/Kg values.
$ws.Range("D141").Value = 44448
$ws.Range("M141").Value = 180
$ws.Range("N141").Value = 21000
$ws.Range("O141").Value = 21000
$ws.Range("P141").Value = 21000
$ws.Range("S141").Value = 1500

# The newly-inserted row 142 holds the data that row 141 used to carry
# before the update above (same market/product info, original date and
# figures).
$ws.Range("A142").Value = 5
$ws.Range("B142").Value = "Macroferia Regional de Talca"
$ws.Range("C142").Value = "Maule"
$ws.Range("D142").Value = 44238
$ws.Range("E142").Value = 7
$ws.Range("F142").Value = "Fruta"
$ws.Range("G142").Value = 100108
$ws.Range("H142").Value = "Tropicales y subtropicales"
$ws.Range("I142").Value = 100108005
$ws.Range("J142").Value = "Piña"
$ws.Range("K142").Value = "Caramelo"
$ws.Range("L142").Value = "Segunda"
$ws.Range("M142").Value = 400
$ws.Range("N142").Value = 12000
$ws.Range("O142").Value = 12000
$ws.Range("P142").Value = 12000
$ws.Range("Q142").Value = "$/caja 14 unidades"
$ws.Range("R142").Value = "Ecuador"
$ws.Range("S142").Value = 857
$ws.Range("T142").Value = 14
